# [EGSVC-36] Added employee details folder
# Rework the "employeeDetails" sheet of eisTestData.xlsx:
#  - drop the EmployeeName / EmployeeCode columns
#  - rename MobileNumber -> Mobile
#  - store DateOfBirth / Mobile / PinCode / DateOfAppointment as plain text
#    instead of numbers / date-serials
#  - add a trailing space to the PermanentAddress string
#  - refresh the view (active sheet / selection / column widths)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("employeeDetails")

# ---------------------------------------------------------------------
# 1. Drop the two leading columns (EmployeeName, EmployeeCode) - this
#    shifts everything else two columns to the left automatically.
# ---------------------------------------------------------------------
$ws1.Columns.Item(2).Delete()
$ws1.Columns.Item(2).Delete()

# ---------------------------------------------------------------------
# 2. Header row
# ---------------------------------------------------------------------
$header = @("dataRow","EmployeeType","Status","DateOfBirth","Gender","MaritalStatus","UserName","IsUserActive","Mobile","PermanentAddress","City","PinCode","DateOfAppointment")
for ($i = 0; $i -lt $header.Length; $i++) {
    $ws1.Cells.Item(1, $i + 1).NumberFormat = "General"
    $ws1.Cells.Item(1, $i + 1).Value = $header[$i]
}

$address = "Municipal Office Rd, N.R.Peta, Near Appollo Hospital, Kurnool, Andhra Pradesh "

# ---------------------------------------------------------------------
# 3. Row 2 (employee1)
# ---------------------------------------------------------------------
$ws1.Cells.Item(2, 1).Value = "employee1"
$ws1.Cells.Item(2, 2).Value = "Permanent"
$ws1.Cells.Item(2, 3).Value = "EMPLOYED"

$ws1.Cells.Item(2, 4).NumberFormat = "@"
$ws1.Cells.Item(2, 4).Value = "01/01/1990"

$ws1.Cells.Item(2, 5).Value = "Male"
$ws1.Cells.Item(2, 6).Value = "UNMARRIED"
$ws1.Cells.Item(2, 7).Value = "testUser1"
$ws1.Cells.Item(2, 8).Value = "Yes"

$ws1.Cells.Item(2, 9).NumberFormat = "@"
$ws1.Cells.Item(2, 9).Value = "9999999999"

$ws1.Cells.Item(2, 10).Value = $address
$ws1.Cells.Item(2, 11).Value = "Kurnool"

$ws1.Cells.Item(2, 12).NumberFormat = "@"
$ws1.Cells.Item(2, 12).Value = "518004"

$ws1.Cells.Item(2, 13).NumberFormat = "@"
$ws1.Cells.Item(2, 13).Value = "01/01/2012"

# ---------------------------------------------------------------------
# 4. Row 3 (employee2)
# ---------------------------------------------------------------------
$ws1.Cells.Item(3, 1).Value = "employee2"
$ws1.Cells.Item(3, 2).Value = "Permanent"
$ws1.Cells.Item(3, 3).Value = "RETIRED"

$ws1.Cells.Item(3, 4).NumberFormat = "@"
$ws1.Cells.Item(3, 4).Value = "02/01/1990"

$ws1.Cells.Item(3, 5).Value = "Female"
$ws1.Cells.Item(3, 6).Value = "MARRIED"
$ws1.Cells.Item(3, 7).Value = "testUser2"
$ws1.Cells.Item(3, 8).Value = "No"

$ws1.Cells.Item(3, 9).NumberFormat = "@"
$ws1.Cells.Item(3, 9).Value = "8888888888"

$ws1.Cells.Item(3, 10).Value = $address
$ws1.Cells.Item(3, 11).Value = "Kurnool"

$ws1.Cells.Item(3, 12).NumberFormat = "@"
$ws1.Cells.Item(3, 12).Value = "518004"

$ws1.Cells.Item(3, 13).NumberFormat = "@"
$ws1.Cells.Item(3, 13).Value = "02/01/2012"

# ---------------------------------------------------------------------
# 5. Column widths (converted from LibreOffice "character width" units
#    stored in the XML to the Excel ColumnWidth property: xml = cw + 5/6)
# ---------------------------------------------------------------------
$widths = @(13.0561224489796,18.0612244897959,15.5612244897959,14.030612244898,9.86224489795918,12.6887755102041,18.0612244897959,16.6683673469388,18.7551020408163,18.6122448979592,13.75,13.8877551020408,21.3928571428571)
for ($i = 0; $i -lt $widths.Length; $i++) {
    $ws1.Columns.Item($i + 1).ColumnWidth = $widths[$i] - 0.8333333333333334
}

# ---------------------------------------------------------------------
# 6. View refresh: employeeDetails becomes the active / selected sheet
#    again, with M3 selected (matches the authored selection).
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("M3").Select()
